# Scheduled market-data refresh: overwrite profit-calculation columns (H:N)
# on specific Leve rows across the job/Sheets, per the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 41676364
$ws.Range("I33").Value = 66668616
$ws.Range("J33").Value = 22610.889
$ws.Range("K33").Value = 66668616
$ws.Range("L33").Value = 22610.889
$ws.Range("M33").Value = -66668387
$ws.Range("N33").Value = -23068.889

# Row 99
$ws.Range("H99").Value = 871.5714
$ws.Range("I99").Value = 788.2
$ws.Range("J99").Value = 1080
$ws.Range("K99").Value = 2364.6
$ws.Range("L99").Value = 3240
$ws.Range("M99").Value = -866.6000000000004
$ws.Range("N99").Value = -6236

# Row 100
$ws.Range("H100").Value = 37039420
$ws.Range("I100").Value = 44446456
$ws.Range("J100").Value = 4250
$ws.Range("K100").Value = 44446456
$ws.Range("L100").Value = 4250
$ws.Range("M100").Value = -44445915
$ws.Range("N100").Value = -5332

# Row 107
$ws.Range("H107").Value = 551.46155
$ws.Range("I107").Value = 403.54544
$ws.Range("J107").Value = 1365
$ws.Range("K107").Value = 403.54544
$ws.Range("L107").Value = 1365
$ws.Range("M107").Value = 1516.45456
$ws.Range("N107").Value = -5205

# Row 132
$ws.Range("H132").Value = 1541918.5
$ws.Range("I132").Value = 1684297.9
$ws.Range("J132").Value = 4221.2
$ws.Range("K132").Value = 5052893.699999999
$ws.Range("L132").Value = 12663.6
$ws.Range("M132").Value = -5050363.699999999
$ws.Range("N132").Value = -17723.6

# Row 138
$ws.Range("H138").Value = 16669771
$ws.Range("I138").Value = 27778568
$ws.Range("J138").Value = 6575.6665
$ws.Range("K138").Value = 83335704
$ws.Range("L138").Value = 19726.9995
$ws.Range("M138").Value = -83330564
$ws.Range("N138").Value = -30006.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 34
$ws.Range("H34").Value = 9798
$ws.Range("J34").Value = 9798
$ws.Range("L34").Value = 9798
$ws.Range("N34").Value = -10340

# Row 61
$ws.Range("H61").Value = 1535.9412
$ws.Range("I61").Value = 1259.25
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 1259.25
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -1047.25
$ws.Range("N61").Value = -2624

# Row 74
$ws.Range("H74").Value = 1252.2858
$ws.Range("I74").Value = 1281.4375
$ws.Range("J74").Value = 1159
$ws.Range("K74").Value = 1281.4375
$ws.Range("L74").Value = 1159
$ws.Range("M74").Value = -407.4375
$ws.Range("N74").Value = -2907

# Row 77
$ws.Range("H77").Value = 1252.2858
$ws.Range("I77").Value = 1281.4375
$ws.Range("J77").Value = 1159
$ws.Range("K77").Value = 6407.1875
$ws.Range("L77").Value = 5795
$ws.Range("M77").Value = -2039.1875
$ws.Range("N77").Value = -14531

# Row 97
$ws.Range("H97").Value = 656.43475
$ws.Range("I97").Value = 656.43475
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 656.43475
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -160.43475
$ws.Range("N97").Value = ""  # cell removed (column now blank for this row)

# Row 132
$ws.Range("H132").Value = 1069.5349
$ws.Range("I132").Value = 1110.7567
$ws.Range("J132").Value = 815.3333
$ws.Range("K132").Value = 3332.2701
$ws.Range("L132").Value = 2445.9999
$ws.Range("M132").Value = -802.2700999999997
$ws.Range("N132").Value = -7505.9999

# Row 136
$ws.Range("H136").Value = 1535.9412
$ws.Range("I136").Value = 1259.25
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 3777.75
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -1227.75
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("BSM")
# Row 21
$ws.Range("H21").Value = 63900
$ws.Range("J21").Value = 63900
$ws.Range("L21").Value = 63900
$ws.Range("N21").Value = -64372

# Row 93
$ws.Range("H93").Value = 37482.668
$ws.Range("J93").Value = 37482.668
$ws.Range("L93").Value = 37482.668
$ws.Range("N93").Value = -41226.668

# Row 94
$ws.Range("H94").Value = 1874.875
$ws.Range("I94").Value = 1742.7142
$ws.Range("J94").Value = 2800
$ws.Range("K94").Value = 1742.7142
$ws.Range("L94").Value = 2800
$ws.Range("M94").Value = -1291.7142
$ws.Range("N94").Value = -3702

# Row 99
$ws.Range("H99").Value = 1635.2632
$ws.Range("I99").Value = 1504.1177
$ws.Range("J99").Value = 2750
$ws.Range("K99").Value = 1504.1177
$ws.Range("L99").Value = 2750
$ws.Range("M99").Value = -6.117700000000013
$ws.Range("N99").Value = -5746

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1773.7216
$ws.Range("I31").Value = 1047.8197
$ws.Range("J31").Value = 4233.722
$ws.Range("K31").Value = 1047.8197
$ws.Range("L31").Value = 4233.722
$ws.Range("M31").Value = -752.8197
$ws.Range("N31").Value = -4823.722

# Row 34
$ws.Range("H34").Value = 1773.7216
$ws.Range("I34").Value = 1047.8197
$ws.Range("J34").Value = 4233.722
$ws.Range("K34").Value = 1047.8197
$ws.Range("L34").Value = 4233.722
$ws.Range("M34").Value = -845.8197
$ws.Range("N34").Value = -4637.722

# Row 58
$ws.Range("H58").Value = 732.34424
$ws.Range("I58").Value = 587.4039
$ws.Range("J58").Value = 1569.7778
$ws.Range("K58").Value = 587.4039
$ws.Range("L58").Value = 1569.7778
$ws.Range("M58").Value = -384.4039
$ws.Range("N58").Value = -1975.7778

# Row 107
$ws.Range("H107").Value = 1608.1818
$ws.Range("I107").Value = 527.1429000000001
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 527.1429000000001
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = 1392.8571
$ws.Range("N107").Value = -7340

# Row 132
$ws.Range("H132").Value = 1195.7544
$ws.Range("I132").Value = 1050.921
$ws.Range("J132").Value = 1485.421
$ws.Range("K132").Value = 3152.763
$ws.Range("L132").Value = 4456.263
$ws.Range("M132").Value = -622.7629999999999
$ws.Range("N132").Value = -9516.262999999999

# Row 134
$ws.Range("H134").Value = 721.0923
$ws.Range("I134").Value = 580.193
$ws.Range("J134").Value = 1725
$ws.Range("K134").Value = 1740.579
$ws.Range("L134").Value = 5175
$ws.Range("M134").Value = 794.421
$ws.Range("N134").Value = -10245

# Row 136
$ws.Range("H136").Value = 732.34424
$ws.Range("I136").Value = 587.4039
$ws.Range("J136").Value = 1569.7778
$ws.Range("K136").Value = 1762.2117
$ws.Range("L136").Value = 4709.3334
$ws.Range("M136").Value = 787.7882999999999
$ws.Range("N136").Value = -9809.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 804.8889
$ws.Range("I129").Value = 604.9167
$ws.Range("J129").Value = 1204.8334
$ws.Range("K129").Value = 1814.7501
$ws.Range("L129").Value = 3614.5002
$ws.Range("M129").Value = 3185.2499
$ws.Range("N129").Value = -13614.5002

# Row 130
$ws.Range("H130").Value = 2496.375
$ws.Range("I130").Value = 1354.2
$ws.Range("J130").Value = 4400
$ws.Range("K130").Value = 4062.6
$ws.Range("L130").Value = 13200
$ws.Range("M130").Value = 957.3999999999996
$ws.Range("N130").Value = -23240

# Row 131
$ws.Range("H131").Value = 915.6087
$ws.Range("I131").Value = 498.6
$ws.Range("K131").Value = 1495.8
$ws.Range("M131").Value = 3544.2

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 890
$ws.Range("I97").Value = 890
$ws.Range("K97").Value = 890
$ws.Range("M97").Value = -394

# Row 122
$ws.Range("H122").Value = 1866.5769
$ws.Range("I122").Value = 1875.9474
$ws.Range("J122").Value = 1841.1428
$ws.Range("K122").Value = 5627.8422
$ws.Range("L122").Value = 5523.428400000001
$ws.Range("M122").Value = -3177.8422
$ws.Range("N122").Value = -10423.4284

# Row 132
$ws.Range("H132").Value = 6154.7393
$ws.Range("I132").Value = 6454.524
$ws.Range("K132").Value = 19363.572
$ws.Range("M132").Value = -16833.572

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5260
$ws.Range("I7").Value = 5750
$ws.Range("J7").Value = 4933.3335
$ws.Range("K7").Value = 5750
$ws.Range("L7").Value = 4933.3335
$ws.Range("M7").Value = -5638
$ws.Range("N7").Value = -5157.3335

# Row 93
$ws.Range("H93").Value = 12229.6
$ws.Range("I93").Value = 21240.4
$ws.Range("J93").Value = 3218.8
$ws.Range("K93").Value = 21240.4
$ws.Range("L93").Value = 3218.8
$ws.Range("M93").Value = -19992.4
$ws.Range("N93").Value = -5714.8

# Row 122
$ws.Range("H122").Value = 3261.8823
$ws.Range("I122").Value = 3080.1333
$ws.Range("J122").Value = 4625
$ws.Range("K122").Value = 9240.3999
$ws.Range("L122").Value = 13875
$ws.Range("M122").Value = -6790.3999
$ws.Range("N122").Value = -18775

# Row 126
$ws.Range("H126").Value = 5260
$ws.Range("I126").Value = 5750
$ws.Range("J126").Value = 4933.3335
$ws.Range("K126").Value = 17250
$ws.Range("L126").Value = 14800.0005
$ws.Range("M126").Value = -14780
$ws.Range("N126").Value = -19740.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2361.2
$ws.Range("I96").Value = 1868.6666
$ws.Range("J96").Value = 3100
$ws.Range("K96").Value = 1868.6666
$ws.Range("L96").Value = 3100
$ws.Range("M96").Value = -495.6666
$ws.Range("N96").Value = -5846

# Row 107
$ws.Range("H107").Value = 4380.577
$ws.Range("I107").Value = 276.3846
$ws.Range("J107").Value = 8484.77
$ws.Range("K107").Value = 829.1537999999999
$ws.Range("L107").Value = 25454.31
$ws.Range("M107").Value = 1090.8462
$ws.Range("N107").Value = -29294.31

